# Estado de Cuenta - update worker list and totals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Prepare row formats before we touch values, so the "plain" (non-bordered)
#    style can be copied from row 16 and the "bottom border" style from the
#    current last data row (22) to the rows that will need them once the
#    table grows from 7 to 10 data rows (16-22 -> 16-25).
# ---------------------------------------------------------------------------

# Copy the bottom-border style (currently on row 22) down to the new last row (25)
$ws.Range("B22:J22").Copy()
$ws.Range("B25:J25").PasteSpecial(-4122)

# Copy the plain row style (row 16) onto the rows that need it:
#  - row 17 is a brand new row
#  - row 22 changes from "last row" (bordered) to a regular row
#  - rows 23 and 24 are brand new rows
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)
$ws.Range("B22:J22").PasteSpecial(-4122)
$ws.Range("B23:J23").PasteSpecial(-4122)
$ws.Range("B24:J24").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Write the updated worker table (rows 16-25)
# ---------------------------------------------------------------------------

$ws.Range("B16").Value2 = "CC"
$ws.Range("C16").Value2 = "45546964"
$ws.Range("D16").Value2 = "ADRIANA LISNEY PEREZ MADIEDO"
$ws.Range("E16").Value2 = "2207"
$ws.Range("F16").Value2 = 60000
$ws.Range("G16").Value2 = 1500000

$ws.Range("B17").Value2 = "CC"
$ws.Range("C17").Value2 = "45546964"
$ws.Range("D17").Value2 = "ADRIANA LISNEY PEREZ MADIEDO"
$ws.Range("E17").Value2 = "2207"
$ws.Range("F17").Value2 = 2000
$ws.Range("G17").Value2 = 1500000

$ws.Range("B18").Value2 = "CC"
$ws.Range("C18").Value2 = "9100670"
$ws.Range("D18").Value2 = "CARLOS AUGUSTO CARDONA RESTREPO"
$ws.Range("E18").Value2 = "1910"
$ws.Range("F18").Value2 = 42000
$ws.Range("G18").Value2 = 952968

$ws.Range("B19").Value2 = "CC"
$ws.Range("C19").Value2 = "73571489"
$ws.Range("D19").Value2 = "LENIN JACOB IBAÑEZ PEREZ"
$ws.Range("E19").Value2 = "2103"
$ws.Range("F19").Value2 = 3511
$ws.Range("G19").Value2 = 908526

$ws.Range("B20").Value2 = "CC"
$ws.Range("C20").Value2 = "73577260"
$ws.Range("D20").Value2 = "MARIO RAFAEL GREY RODRIGUEZ"
$ws.Range("E20").Value2 = "2109"
$ws.Range("F20").Value2 = 1211
$ws.Range("G20").Value2 = 908526

$ws.Range("B21").Value2 = "CC"
$ws.Range("C21").Value2 = "9145170"
$ws.Range("D21").Value2 = "HUGO ALFONSO MERCADO ZABALETA"
$ws.Range("E21").Value2 = "2201"
$ws.Range("F21").Value2 = 47137
$ws.Range("G21").Value2 = 1178421

$ws.Range("B22").Value2 = "CC"
$ws.Range("C22").Value2 = "73122996"
$ws.Range("D22").Value2 = "ALEJANDRO ISMAEL GUETTE SAAVEDRA"
$ws.Range("E22").Value2 = "2110"
$ws.Range("F22").Value2 = 36341
$ws.Range("G22").Value2 = 908526

$ws.Range("B23").Value2 = "CC"
$ws.Range("C23").Value2 = "73350017"
$ws.Range("D23").Value2 = "RICARDO POLO PATERNINA"
$ws.Range("E23").Value2 = "1911"
$ws.Range("F23").Value2 = 33125
$ws.Range("G23").Value2 = 877803

$ws.Range("B24").Value2 = "CC"
$ws.Range("C24").Value2 = "1047503421"
$ws.Range("D24").Value2 = "LUIS ALEJANDRO PUERTA DOMINGUEZ"
$ws.Range("E24").Value2 = "2507"
$ws.Range("F24").Value2 = 56940
$ws.Range("G24").Value2 = 877803

$ws.Range("B25").Value2 = "CC"
$ws.Range("C25").Value2 = "1052075224"
$ws.Range("D25").Value2 = "WILMER ALBERTO VARGAS ZAPATA"
$ws.Range("E25").Value2 = "2103"
$ws.Range("F25").Value2 = 68000
$ws.Range("G25").Value2 = 1700000

# ---------------------------------------------------------------------------
# 3) Update the summary header values (Valor Mora, Cant. Trabajadores, Cant. Periodos)
# ---------------------------------------------------------------------------

$ws.Range("E11").Value2 = 350265
$ws.Range("C13").Value2 = 9
$ws.Range("F13").Value2 = 8

# ---------------------------------------------------------------------------
# 4) Move the signature footer block down from rows 27-28 to rows 30-31
#    (table grew by 3 rows so the footer has to move down to stay below it)
# ---------------------------------------------------------------------------

$ws.Range("B27:C27").Copy()
$ws.Range("B30:C30").PasteSpecial(-4122)
$ws.Range("H27:J27").Copy()
$ws.Range("H30:J30").PasteSpecial(-4122)

$ws.Range("B28:C28").Copy()
$ws.Range("B31:C31").PasteSpecial(-4122)
$ws.Range("H28:J28").Copy()
$ws.Range("H31:J31").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$ws.Range("B27:J28").Clear()

$ws.Range("B27:C27").UnMerge()
$ws.Range("B28:C28").UnMerge()
$ws.Range("H27:J27").UnMerge()
$ws.Range("H28:J28").UnMerge()

$ws.Range("B30").Value2 = "___________________________________"
$ws.Range("H30").Value2 = "___________________________________"
$ws.Range("B31").Value2 = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Range("H31").Value2 = "FIRMA DEL REPRESENTANTE LEGAL"

$ws.Range("B30:C30").Merge()
$ws.Range("B31:C31").Merge()
$ws.Range("H30:J30").Merge()
$ws.Range("H31:J31").Merge()

$wb.Save()
